$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 45,4
$arr[0,0] = "2024-12-19 00:28:08"
$arr[0,1] = -0.1215316019114107
$arr[0,2] = -0.001756474862245996
$arr[0,3] = 0.008538688149035216
$arr[1,0] = "2024-12-19 00:28:09"
$arr[1,1] = -0.1207703233562643
$arr[1,2] = -0.001735773831951996
$arr[1,3] = 0.00838519867832738
$arr[2,0] = "2024-12-19 00:28:10"
$arr[2,1] = -0.1226051998737967
$arr[2,2] = -0.001938380737249996
$arr[2,3] = 0.009506222308882123
$arr[3,0] = "2024-12-19 00:28:11"
$arr[3,1] = -0.121655228343443
$arr[3,2] = -0.001829814209179995
$arr[3,3] = 0.008904258617754757
$arr[4,0] = "2024-12-19 00:28:12"
$arr[4,1] = -0.1236755445090239
$arr[4,2] = -0.001855829684903995
$arr[4,3] = 0.009180829871860474
$arr[5,0] = "2024-12-19 00:28:13"
$arr[5,1] = -0.1240561837865971
$arr[5,2] = -0.001928612280411996
$arr[5,3] = 0.009570251180475149
$arr[6,0] = "2024-12-19 00:28:14"
$arr[6,1] = -0.1220326142938575
$arr[6,2] = -0.001890196432017995
$arr[6,3] = 0.009226624485123102
$arr[7,0] = "2024-12-19 00:28:15"
$arr[7,1] = -0.1194722458199249
$arr[7,2] = -0.001715376484253996
$arr[7,3] = 0.008197595240020474
$arr[8,0] = "2024-12-19 00:28:16"
$arr[8,1] = -0.1217007749236654
$arr[8,2] = -0.001633280955801996
$arr[8,3] = 0.007950862319566716
$arr[9,0] = "2024-12-19 00:28:17"
$arr[9,1] = -0.1229728258427349
$arr[9,2] = -0.001750451824091995
$arr[9,3] = 0.008610320292406503
$arr[10,0] = "2024-12-19 00:28:18"
$arr[10,1] = -0.122357947009732
$arr[10,2] = -0.001991120281421995
$arr[10,3] = 0.009745175595369408
$arr[11,0] = "2024-12-19 00:28:19"
$arr[11,1] = -0.124658049310965
$arr[11,2] = -0.002034850575245996
$arr[11,3] = 0.01014642013397843
$arr[12,0] = "2024-12-19 00:28:20"
$arr[12,1] = -0.1213884555164259
$arr[12,2] = -0.001868533740169996
$arr[12,3] = 0.009072736991982659
$arr[13,0] = "2024-12-19 00:28:21"
$arr[13,1] = -0.1213949621707434
$arr[13,2] = -0.001837608729143996
$arr[13,3] = 0.008923057686362527
$arr[14,0] = "2024-12-19 00:28:22"
$arr[14,1] = -0.1206239236341207
$arr[14,2] = -0.001856993801521996
$arr[14,3] = 0.008959915140152992
$arr[15,0] = "2024-12-19 00:28:23"
$arr[15,1] = -0.1195958722519572
$arr[15,2] = -0.001927498777559996
$arr[15,3] = 0.009220835902674756
$arr[16,0] = "2024-12-19 00:28:24"
$arr[16,1] = -0.1201164045973564
$arr[16,2] = -0.001845099566511996
$arr[16,3] = 0.008865069042142473
$arr[17,0] = "2024-12-19 00:28:25"
$arr[17,1] = -0.1193876593137975
$arr[17,2] = -0.002009341237181995
$arr[17,3] = 0.009595621882793934
$arr[18,0] = "2024-12-19 00:28:26"
$arr[18,1] = -0.1216877616150305
$arr[18,2] = -0.001776517913581995
$arr[18,3] = 0.008647219534911887
$arr[19,0] = "2024-12-19 00:28:27"
$arr[19,1] = -0.1218536813001265
$arr[19,2] = -0.001812048777313996
$arr[19,3] = 0.00883219256844414
$arr[20,0] = "2024-12-19 00:28:28"
$arr[20,1] = -0.122410000244272
$arr[20,2] = -0.001748781569813996
$arr[20,3] = 0.008562734095524379
$arr[21,0] = "2024-12-19 00:28:29"
$arr[21,1] = -0.1225206133676693
$arr[21,2] = -0.001802887685667996
$arr[21,3] = 0.008835636203242426
$arr[22,0] = "2024-12-19 00:28:30"
$arr[22,1] = -0.120952509677154
$arr[22,2] = -0.002048769360895997
$arr[22,3] = 0.009912151838001187
$arr[23,0] = "2024-12-19 00:28:31"
$arr[23,1] = -0.1210663761277101
$arr[23,2] = -0.001813061052633996
$arr[23,3] = 0.008780029253627572
$arr[24,0] = "2024-12-19 00:28:32"
$arr[24,1] = -0.1201164045973564
$arr[24,2] = -0.001892777734083995
$arr[24,3] = 0.009094146244804029
$arr[25,0] = "2024-12-19 00:28:33"
$arr[25,1] = -0.1186686740117148
$arr[25,2] = -0.001811390798355996
$arr[25,3] = 0.008598213766317099
$arr[26,0] = "2024-12-19 00:28:34"
$arr[26,1] = -0.1205328304736759
$arr[26,2] = -0.001889133542931996
$arr[26,3] = 0.009108104522894275
$arr[27,0] = "2024-12-19 00:28:35"
$arr[27,1] = -0.1200350714183878
$arr[27,2] = -0.001836039702397995
$arr[27,3] = 0.008815566272173555
$arr[28,0] = "2024-12-19 00:28:36"
$arr[28,1] = -0.1209492563499953
$arr[28,2] = -0.001742303007765996
$arr[28,3] = 0.008429210125026289
$arr[29,0] = "2024-12-19 00:28:37"
$arr[29,1] = -0.1187304872277309
$arr[29,2] = -0.001708543625843996
$arr[29,3] = 0.008114248685851665
$arr[30,0] = "2024-12-19 00:28:38"
$arr[30,1] = -0.1189321935115731
$arr[30,2] = -0.001730712455351995
$arr[30,3] = 0.008233497146112535
$arr[31,0] = "2024-12-19 00:28:39"
$arr[31,1] = -0.1187239805734134
$arr[31,2] = -0.001820804958831995
$arr[31,3] = 0.008646928502413788
$arr[32,0] = "2024-12-19 00:28:40"
$arr[32,1] = -0.1182880347341415
$arr[32,2] = -0.001860891061503996
$arr[32,3] = 0.008804845860785527
$arr[33,0] = "2024-12-19 00:28:41"
$arr[33,1] = -0.1185092609809362
$arr[33,2] = -0.001868331285105996
$arr[33,3] = 0.008856582394618978
$arr[34,0] = "2024-12-19 00:28:42"
$arr[34,1] = -0.1202432843565475
$arr[34,2] = -0.001830978325797996
$arr[34,3] = 0.008806513899184147
$arr[35,0] = "2024-12-19 00:28:43"
$arr[35,1] = -0.1196251521963859
$arr[35,2] = -0.001718565151511995
$arr[35,3] = 0.008223344712361099
$arr[36,0] = "2024-12-19 00:28:44"
$arr[36,1] = -0.1193356060792576
$arr[36,2] = -0.001537266641699995
$arr[36,3] = 0.007338025855707753
$arr[37,0] = "2024-12-19 00:28:45"
$arr[37,1] = -0.1203538974799449
$arr[37,2] = -0.001620323831705996
$arr[37,3] = 0.007800491533018197
$arr[38,0] = "2024-12-19 00:28:46"
$arr[38,1] = -0.11852552761673
$arr[38,2] = -0.001698927010303996
$arr[38,3] = 0.008054648811143789
$arr[39,0] = "2024-12-19 00:28:47"
$arr[39,1] = -0.1203181108811986
$arr[39,2] = -0.001665370083445997
$arr[39,3] = 0.00801496729433146
$arr[40,0] = "2024-12-19 00:28:48"
$arr[40,1] = -0.1208711764981854
$arr[40,2] = -0.001982161644839995
$arr[40,3] = 0.009583448400855541
$arr[41,0] = "2024-12-19 00:28:49"
$arr[41,1] = -0.1213754422077909
$arr[41,2] = -0.001828397023731996
$arr[41,3] = 0.008876899891475193
$arr[42,0] = "2024-12-19 00:28:50"
$arr[42,1] = -0.119943978257943
$arr[42,2] = -0.001981048141987996
$arr[42,3] = 0.009504591810821859
$arr[43,0] = "2024-12-19 00:28:51"
$arr[43,1] = -0.1197813119000057
$arr[43,2] = -0.001961410000779997
$arr[43,3] = 0.009397610522688764
$arr[44,0] = "2024-12-19 00:28:52"
$arr[44,1] = -0.1197845652271644
$arr[44,2] = -0.002017641894805995
$arr[44,3] = 0.009667294286137934

$range = $ws.Range("A204:D248")
$range.Value = $arr

Write-Host "Added rows from row" 204 "to row" 248
